# "Data Driven from excel sheet"
# Add a second sheet ("InvalidLogin") next to the existing "ValidLogin" sheet,
# holding a UserName/Password pair used for a negative login test case.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

# The workbook currently has a single sheet ("ValidLogin") - insert the new
# sheet right after it so it lands as the second tab.
$validLoginSheet = $sheets.Item(1)
$invalidLoginSheet = $sheets.Add([System.Reflection.Missing]::Value, $validLoginSheet)
$invalidLoginSheet.Name = "InvalidLogin"

$invalidLoginSheet.Range("A1").Value = "UserName"
$invalidLoginSheet.Range("B1").Value = "Password"
$invalidLoginSheet.Range("A2").Value = "asd"
$invalidLoginSheet.Range("B2").Value = "dd"

# Leave the same cell selected as in the authored workbook.
$invalidLoginSheet.Range("B2").Select() | Out-Null
